$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Copy the date number-format from an existing filled-in "Date" cell (B25)
# onto the three new rows being populated (B26, B27, B28) so they render as
# dates like the rest of the log, then set their values.
# ---------------------------------------------------------------------------
$ws.Range("B25").Copy() | Out-Null
$ws.Range("B26:B28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row 26 - Entry 10
# ---------------------------------------------------------------------------
$ws.Range("B26").Value = 46062
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = "OpenCodeAI"
$ws.Range("E26").Value = "Asked OpenCode to just add whatever functions it wants to add that it thinks a user that has logged in should do"
$ws.Range("F26").Value = "The AI added loads of features, User profile, Filter System, Reviews, rating system, favourites, pre ordering concessions, promotional codes, booking history and even made it so that the user data gets saved in localStorage"
$ws.Range("G26").Value = "The user dashboard worked, comments worked partially, concessions ordering worked fine and viewing bookings that had been made was fine"
$ws.Range("H26").Value = "Users were still able to make multiple bookings for the same seats, comments couldnt be deleted. only updated"
$ws.Range("I26").Value = "I went through all the features, seeing if i could comment twice, and checked if the bookings changed anything in the seating once had booked"
$ws.Range("J26").Value = "might be a few features i forgot to check or glossed over as im not completely sure what the ai added"

# ---------------------------------------------------------------------------
# Row 27 - Entry 11
# ---------------------------------------------------------------------------
$ws.Range("B27").Value = 46062
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = "OpenCodeAI"
$ws.Range("E27").Value = "Got openCode to remake the website's design, and add more cinema data"
$ws.Range("F27").Value = "The ai added filtering to the map, redesigned the page entirely."
$ws.Range("G27").Value = "browsing the website worked but thats about it"
$ws.Range("H27").Value = "The search function didnt work at all nor the filtering, a notifcation would just come up saying what it was filtering for, it also deleted the images i had put in for copilot to use whenever i clicked on details on the map it would also put my user name in the search bar. when the ai was trying to create info for one of the movies there were multiple directors and the ai didnt but them in an array"
$ws.Range("I27").Value = "went through the entire page and carefully checked every detail that the ai had said it had changed"

# ---------------------------------------------------------------------------
# Row 28 - Entry 12
# ---------------------------------------------------------------------------
$ws.Range("B28").Value = 46062
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = "OpenCodeAI"
$ws.Range("E28").Value = "Got the ai to fix whatever errors it had and to remove any buttons that were useless"
$ws.Range("F28").Value = "the ai removed most of the buttons and completely fixed the search function and made it so that you can search by director"
$ws.Range("G28").Value = "practically everything worked perfectly fine to the point that even made it so that the user can search by director and rating"
$ws.Range("H28").Value = "The map is still a bit broken when clicking on details for the cinema it will just zoom in on the point on the map"
$ws.Range("I28").Value = "Messed with the search engine to see what details the user can search with and i rechecked the map "
$ws.Range("J28").Value = "Might've missed a button or two that the ai had put in that i hadnt noticed"

# ---------------------------------------------------------------------------
# Update the saved view/selection to match where the author left off editing
# ---------------------------------------------------------------------------
$window = $excel.ActiveWindow
$window.ScrollRow = 18
$ws.Range("K27").Select() | Out-Null
